# "removed distinction between argument and statement"
#
# Merge the "Arguments" sheet's legalGround / supports / evidence columns
# into the "Statements" sheet (which has a spare, always-empty "Argument"
# column), then delete the now-redundant "Arguments" sheet.

$wb = $excel.ActiveWorkbook

$wsStatements = $wb.Worksheets.Item("Statements")
$wsArguments  = $wb.Worksheets.Item("Arguments")

# Move a single cell's value+formatting from src to dst. When the source
# cell is blank, just clear the destination instead of copying -- a
# straight Range.Copy of a blank cell does not clear an already-populated
# destination cell, and we don't want to leave stray empty cell nodes
# behind either.
function Move-Cell($srcWs, $srcAddr, $dstWs, $dstAddr) {
    $src = $srcWs.Range($srcAddr)
    $dst = $dstWs.Range($dstAddr)
    if ($src.Value2 -eq $null) {
        $dst.ClearContents()
    } else {
        $src.Copy($dst)
    }
}

# --- 1. Shift the Statements data columns left ---------------------------
# Before: A=[Statement] B=Argument(unused) C=madeBy D=moment  E=phrasing
# After : A=[Statement] B=madeBy           C=moment D=phrasing
# Go column-by-column, left to right, so each source column is fully read
# before it later becomes a copy destination.
foreach ($r in 1..8) {
    $srcAddr = "C" + $r
    $dstAddr = "B" + $r
    Move-Cell $wsStatements $srcAddr $wsStatements $dstAddr
}
foreach ($r in 1..8) {
    $srcAddr = "D" + $r
    $dstAddr = "C" + $r
    Move-Cell $wsStatements $srcAddr $wsStatements $dstAddr
}
foreach ($r in 1..8) {
    $srcAddr = "E" + $r
    $dstAddr = "D" + $r
    Move-Cell $wsStatements $srcAddr $wsStatements $dstAddr
}

# --- 2. Pull the Arguments columns (legalGround/supports/evidence) in ---
# as the new E:G columns.
#  * Header rows line up 1:1 (row 1 -> row 1, row 2 -> row 2).
#  * Data rows shift down by one: Arguments row 3 (its first data row,
#    argument "s2") lands on Statements row 4 (the "s2" statement row),
#    and the two trailing Arguments rows that have no matching Statement
#    ("s5"/"s6") simply land on new, otherwise-empty rows 9 and 10.
# Row 3 ("s1") has no corresponding Argument row at all, so the leftover
# E3 (old phrasing text, already copied to D3 above) needs to be wiped.
$wsStatements.Range("E3:G3").ClearContents()

$srcCols = @("B", "C", "D")
$dstCols = @("E", "F", "G")

foreach ($r in 1..2) {
    for ($i = 0; $i -lt 3; $i++) {
        $srcAddr = $srcCols[$i] + $r
        $dstAddr = $dstCols[$i] + $r
        Move-Cell $wsArguments $srcAddr $wsStatements $dstAddr
    }
}

foreach ($r in 3..9) {
    $dstR = $r + 1
    for ($i = 0; $i -lt 3; $i++) {
        $srcAddr = $srcCols[$i] + $r
        $dstAddr = $dstCols[$i] + $dstR
        Move-Cell $wsArguments $srcAddr $wsStatements $dstAddr
    }
}

$wsStatements.Range("D15").Select()

# --- 3. Remove the now-redundant Arguments sheet -------------------------
$excel.DisplayAlerts = $false
$wsArguments.Delete()
$excel.DisplayAlerts = $true
